$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(286,1).Value = "x"
$win = $wb.Windows.Item(1)
$win.ScrollRow = 271
$win.ScrollColumn = 1
